$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data values of row 2 and row 3 for columns B, D, E, G, H
# (these are the columns whose values actually differ between the two rows)
$cols = @("B", "D", "E", "G", "H")

foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}
